# april 20 pm refresh
# The "Data providers" sub-heading in the "Data source type and data
# collection method" table was hand-formatted (explicit pBdr/shd/color/
# size/font overrides) instead of using the shared "M.Header" paragraph
# style that every other section heading in the document already uses
# (Institutional information, Concepts and definitions, Calendar, Data
# compilers, Methodology, References, ...). Bring it in line with the
# rest of the document by applying that style, which also carries the
# blue bottom-border / blue-18pt-text look the manual formatting was
# trying to reproduce.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Data providers") {
        $p.Style = "MHeader"
    }
}
